$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from Ferlab.bio CodeS" worksheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- 2. Metadata worksheet updates ---
$ws = $wb.Worksheets.Item(1)

# 2a. Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 2b. Update the Contact value (row 10, column B)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 2c. Insert a new "Jurisdiction" row before the current row 11 ("Description"),
#     pushing the existing rows 11-14 down to 12-15. Capture the values first
#     so the cascading writes below don't clobber data we still need to read.
$a11 = $ws.Range("A11").Text
$b11 = $ws.Range("B11").Text
$a12 = $ws.Range("A12").Text
$b12 = $ws.Range("B12").Text
$a13 = $ws.Range("A13").Text
$b13 = $ws.Range("B13").Text
$a14 = $ws.Range("A14").Text
$b14 = $ws.Range("B14").Text

# Give the brand-new row 15 the same formatting (style) as the row it is
# replacing the content of (row 14), before we overwrite the values.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12
$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

Write-Host "edits applied"
